$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = '5840535 - Messias Borges Silva'
$ws.Range("C10").Value = '5840535 - Messias Borges Silva'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Introduction Conventional Experimentation, Full Factorial Experiments, Fractional Factorial Experiments, Analysis of Variance, Response Surface Methodology, Taguchi’s Method'
$ws.Range("C14").Value = 'Introduction Conventional Experimentation, Full Factorial Experiments, Fractional Factorial Experiments, Analysis of Variance, Response Surface Methodology, Taguchi’s Method'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2018'
$ws.Range("C15").Value = '01/01/2018'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = 'Introduction Conventional Experimentation, Full Factorial Experiments (2k), Fractional Factorial Experiments 2k-p , Plackett Burman Method, Analysis of Variance, Response Surface Methodology, Taguchi Method'
$ws.Range("C16").Value = 'Introduction Conventional Experimentation, Full Factorial Experiments (2k), Fractional Factorial Experiments 2k-p , Plackett Burman Method, Analysis of Variance, Response Surface Methodology, Taguchi Method'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840535 - Messias Borges Silva'
$ws.Range("C18").Value = '5840535 - Messias Borges Silva'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = 'LOB1049 -  Estatística Multivariada  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1049 -  Estatística Multivariada  (Requisito fraco)
'

$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).EntireRow.AutoFit()
$ws.Rows.Item(23).RowHeight = 30

$ws.Rows.Item(24).Delete()

Write-Output "done"
